$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ between row 16 and row 17 (a full
# row swap occurred between these two records in the source data).
$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")

foreach ($col in $cols) {
    $c16 = $ws.Range($col + "16")
    $c17 = $ws.Range($col + "17")
    $v16 = $c16.Value2
    $v17 = $c17.Value2
    $c16.Value = $v17
    $c17.Value = $v16
}
